$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3359.1667
$ws.Cells.Item(40, 9).Value = 5250.5
$ws.Cells.Item(40, 11).Value = 5250.5
$ws.Cells.Item(40, 13).Value = -5075.5
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 13).ClearContents()
$ws.Cells.Item(109, 8).Value = 73170.836
$ws.Cells.Item(109, 9).Value = 45000
$ws.Cells.Item(109, 11).Value = 45000
$ws.Cells.Item(109, 13).Value = -43613
$ws.Cells.Item(112, 8).Value = 6410.6523
$ws.Cells.Item(112, 10).Value = 7394.4736
$ws.Cells.Item(112, 12).Value = 22183.4208
$ws.Cells.Item(112, 14).Value = -24399.4208
$ws.Cells.Item(116, 8).Value = 357037.66
$ws.Cells.Item(116, 10).Value = 454306.5
$ws.Cells.Item(116, 12).Value = 454306.5
$ws.Cells.Item(116, 14).Value = -461190.5
$ws.Cells.Item(137, 8).Value = 2217.8
$ws.Cells.Item(137, 9).Value = 1869.7931
$ws.Cells.Item(137, 11).Value = 5609.379300000001
$ws.Cells.Item(137, 13).Value = -3059.379300000001
$ws.Cells.Item(138, 8).Value = 1721.2979
$ws.Cells.Item(138, 9).Value = 1107.0605
$ws.Cells.Item(138, 10).Value = 3169.1428
$ws.Cells.Item(138, 11).Value = 3321.1815
$ws.Cells.Item(138, 12).Value = 9507.428400000001
$ws.Cells.Item(138, 13).Value = 1818.8185
$ws.Cells.Item(138, 14).Value = -19787.4284
$ws.Cells.Item(141, 8).Value = 729.881
$ws.Cells.Item(141, 9).Value = 735.70734
$ws.Cells.Item(141, 11).Value = 2207.12202
$ws.Cells.Item(141, 13).Value = 2972.87798

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 158
$ws.Cells.Item(5, 9).Value = 158
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 158
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).Value = -46
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 13).ClearContents()
$ws.Cells.Item(11, 8).Value = 1401620
$ws.Cells.Item(11, 9).Value = 2333533.2
$ws.Cells.Item(11, 10).Value = 3750
$ws.Cells.Item(11, 11).Value = 2333533.2
$ws.Cells.Item(11, 12).Value = 3750
$ws.Cells.Item(11, 13).Value = -2333389.2
$ws.Cells.Item(11, 14).Value = -4038
$ws.Cells.Item(74, 8).Value = 6703
$ws.Cells.Item(74, 9).Value = 7090.5835
$ws.Cells.Item(74, 10).Value = 5374.143
$ws.Cells.Item(74, 11).Value = 7090.5835
$ws.Cells.Item(74, 12).Value = 5374.143
$ws.Cells.Item(74, 13).Value = -6216.5835
$ws.Cells.Item(74, 14).Value = -7122.143
$ws.Cells.Item(77, 8).Value = 6703
$ws.Cells.Item(77, 9).Value = 7090.5835
$ws.Cells.Item(77, 10).Value = 5374.143
$ws.Cells.Item(77, 11).Value = 35452.9175
$ws.Cells.Item(77, 12).Value = 26870.715
$ws.Cells.Item(77, 13).Value = -31084.9175
$ws.Cells.Item(77, 14).Value = -35606.715
$ws.Cells.Item(122, 8).Value = 2736.25
$ws.Cells.Item(122, 9).Value = 2648.72
$ws.Cells.Item(122, 10).Value = 3465.6667
$ws.Cells.Item(122, 11).Value = 7946.16
$ws.Cells.Item(122, 12).Value = 10397.0001
$ws.Cells.Item(122, 13).Value = -5496.16
$ws.Cells.Item(122, 14).Value = -15297.0001
$ws.Cells.Item(132, 8).Value = 2386.0833
$ws.Cells.Item(132, 9).Value = 2407.9702
$ws.Cells.Item(132, 11).Value = 7223.910600000001
$ws.Cells.Item(132, 13).Value = -4693.910600000001
$ws.Cells.Item(134, 8).Value = 78636.75
$ws.Cells.Item(134, 10).Value = 78636.75
$ws.Cells.Item(134, 12).Value = 78636.75
$ws.Cells.Item(134, 14).Value = -88776.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 158
$ws.Cells.Item(4, 9).Value = 158
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 158
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).Value = -43
$ws.Cells.Item(105, 8).Value = 3599.7
$ws.Cells.Item(105, 9).Value = 2888.5557
$ws.Cells.Item(105, 11).Value = 2888.5557
$ws.Cells.Item(105, 13).Value = -1141.5557
$ws.Cells.Item(107, 8).Value = 37041160
$ws.Cells.Item(107, 9).Value = 6406.3335
$ws.Cells.Item(107, 11).Value = 6406.3335
$ws.Cells.Item(107, 13).Value = -4486.3335
$ws.Cells.Item(134, 8).Value = 2693.4
$ws.Cells.Item(134, 9).Value = 1602.4783
$ws.Cells.Item(134, 11).Value = 4807.4349
$ws.Cells.Item(134, 13).Value = -2272.4349

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 70000
$ws.Cells.Item(20, 10).Value = 70000
$ws.Cells.Item(20, 12).Value = 70000
$ws.Cells.Item(20, 14).Value = -70472
$ws.Cells.Item(30, 8).Value = 70000
$ws.Cells.Item(30, 10).Value = 70000
$ws.Cells.Item(30, 12).Value = 70000
$ws.Cells.Item(30, 14).Value = -70182
$ws.Cells.Item(31, 8).Value = 85054.12
$ws.Cells.Item(31, 9).Value = 120886.37
$ws.Cells.Item(31, 11).Value = 120886.37
$ws.Cells.Item(31, 13).Value = -120591.37
$ws.Cells.Item(34, 8).Value = 85054.12
$ws.Cells.Item(34, 9).Value = 120886.37
$ws.Cells.Item(34, 11).Value = 120886.37
$ws.Cells.Item(34, 13).Value = -120684.37
$ws.Cells.Item(50, 8).Value = 60000
$ws.Cells.Item(50, 10).Value = 60000
$ws.Cells.Item(50, 12).Value = 60000
$ws.Cells.Item(50, 14).Value = -61250
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).ClearContents()
$ws.Cells.Item(58, 8).Value = 3948.9167
$ws.Cells.Item(58, 9).Value = 1192.8
$ws.Cells.Item(58, 11).Value = 1192.8
$ws.Cells.Item(58, 13).Value = -989.8
$ws.Cells.Item(60, 8).Value = 17496.428
$ws.Cells.Item(60, 10).Value = 25000
$ws.Cells.Item(60, 12).Value = 25000
$ws.Cells.Item(60, 14).Value = -26022
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(122, 8).Value = 2119
$ws.Cells.Item(122, 9).Value = 1415.8334
$ws.Cells.Item(122, 11).Value = 4247.5002
$ws.Cells.Item(122, 13).Value = -1797.5002
$ws.Cells.Item(128, 8).Value = 70000
$ws.Cells.Item(128, 10).Value = 70000
$ws.Cells.Item(128, 12).Value = 70000
$ws.Cells.Item(128, 14).Value = -79960
$ws.Cells.Item(132, 8).Value = 5817473
$ws.Cells.Item(132, 9).Value = 3633.45
$ws.Cells.Item(132, 11).Value = 10900.35
$ws.Cells.Item(132, 13).Value = -8370.349999999999
$ws.Cells.Item(136, 8).Value = 3948.9167
$ws.Cells.Item(136, 9).Value = 1192.8
$ws.Cells.Item(136, 11).Value = 3578.4
$ws.Cells.Item(136, 13).Value = -1028.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 5000525
$ws.Cells.Item(11, 9).Value = 6666666.5
$ws.Cells.Item(11, 10).Value = 2100
$ws.Cells.Item(11, 11).Value = 6666666.5
$ws.Cells.Item(11, 12).Value = 2100
$ws.Cells.Item(11, 13).Value = -6666527.5
$ws.Cells.Item(11, 14).Value = -2378
$ws.Cells.Item(113, 8).Value = 2828.1904
$ws.Cells.Item(113, 9).Value = 2283.2856
$ws.Cells.Item(113, 10).Value = 3918
$ws.Cells.Item(113, 11).Value = 2283.2856
$ws.Cells.Item(113, 12).Value = 3918
$ws.Cells.Item(113, 13).Value = -113.2856000000002
$ws.Cells.Item(113, 14).Value = -8258
$ws.Cells.Item(122, 8).Value = 15628590
$ws.Cells.Item(122, 9).Value = 17243066
$ws.Cells.Item(122, 10).Value = 21999.666
$ws.Cells.Item(122, 11).Value = 51729198
$ws.Cells.Item(122, 12).Value = 65998.99800000001
$ws.Cells.Item(122, 13).Value = -51726748
$ws.Cells.Item(122, 14).Value = -70898.99800000001
$ws.Cells.Item(132, 8).Value = 278157.38
$ws.Cells.Item(132, 9).Value = 469842
$ws.Cells.Item(132, 11).Value = 1409526
$ws.Cells.Item(132, 13).Value = -1406996

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12, 8).Value = 2080
$ws.Cells.Item(12, 9).Value = 1734.2858
$ws.Cells.Item(12, 10).Value = 4500
$ws.Cells.Item(12, 11).Value = 1734.2858
$ws.Cells.Item(12, 12).Value = 4500
$ws.Cells.Item(12, 13).Value = -1564.2858
$ws.Cells.Item(12, 14).Value = -4840
$ws.Cells.Item(45, 8).Value = 4400
$ws.Cells.Item(45, 10).Value = 1666.6666
$ws.Cells.Item(45, 12).Value = 1666.6666
$ws.Cells.Item(45, 14).Value = -2480.6666
$ws.Cells.Item(55, 8).Value = 821.4737
$ws.Cells.Item(55, 9).Value = 263.875
$ws.Cells.Item(55, 10).Value = 1227
$ws.Cells.Item(55, 11).Value = 263.875
$ws.Cells.Item(55, 12).Value = 1227
$ws.Cells.Item(55, 13).Value = -90.875
$ws.Cells.Item(55, 14).Value = -1573
$ws.Cells.Item(87, 8).Value = 60141.75
$ws.Cells.Item(87, 9).Value = 60000
$ws.Cells.Item(87, 11).Value = 60000
$ws.Cells.Item(87, 13).Value = -58877
$ws.Cells.Item(90, 8).Value = 60141.75
$ws.Cells.Item(90, 9).Value = 60000
$ws.Cells.Item(90, 11).Value = 180000
$ws.Cells.Item(90, 13).Value = -174384
$ws.Cells.Item(122, 8).Value = 2749.75
$ws.Cells.Item(122, 9).Value = 2434.3076
$ws.Cells.Item(122, 11).Value = 7302.9228
$ws.Cells.Item(122, 13).Value = -4852.9228
$ws.Cells.Item(127, 8).Value = 165445
$ws.Cells.Item(127, 10).Value = 165445
$ws.Cells.Item(127, 12).Value = 165445
$ws.Cells.Item(127, 14).Value = -175365
$ws.Cells.Item(136, 8).Value = 45098.793
$ws.Cells.Item(136, 9).Value = 2759
$ws.Cells.Item(136, 10).Value = 412043.66
$ws.Cells.Item(136, 11).Value = 8277
$ws.Cells.Item(136, 12).Value = 1236130.98
$ws.Cells.Item(136, 13).Value = -5727
$ws.Cells.Item(136, 14).Value = -1241230.98

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 13002
$ws.Cells.Item(17, 9).Value = 6004
$ws.Cells.Item(17, 10).Value = 20000
$ws.Cells.Item(17, 11).Value = 6004
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = -5832
$ws.Cells.Item(17, 14).Value = -20344
$ws.Cells.Item(122, 8).Value = 2017.6923
$ws.Cells.Item(122, 9).Value = 1930.091
$ws.Cells.Item(122, 10).Value = 2499.5
$ws.Cells.Item(122, 11).Value = 5790.272999999999
$ws.Cells.Item(122, 12).Value = 7498.5
$ws.Cells.Item(122, 13).Value = -3340.272999999999
$ws.Cells.Item(122, 14).Value = -12398.5
$ws.Cells.Item(125, 8).Value = 33333
$ws.Cells.Item(125, 10).Value = 33333
$ws.Cells.Item(125, 12).Value = 33333
$ws.Cells.Item(125, 14).Value = -43173
$ws.Cells.Item(136, 8).Value = 9056797
$ws.Cells.Item(136, 9).Value = 10556140
$ws.Cells.Item(136, 10).Value = 60742.5
$ws.Cells.Item(136, 11).Value = 31668420
$ws.Cells.Item(136, 12).Value = 182227.5
$ws.Cells.Item(136, 13).Value = -31665870
$ws.Cells.Item(136, 14).Value = -187327.5
